$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Day"
$ws.Range("B1").Value = "Money Spent"
$ws.Range("C1").Value = "Description"
